# Update the "Leases" sheet (the active sheet) with new test data and
# adjust the selected cell to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values that changed in the shared strings table:
#   Ratan                 -> Chandra
#   Tata                   -> Siddarth
#   ratan.tata@gmail.com   -> chandra.siddarth@yahoo.com
$ws.Range("A2").Value = "Chandra"
$ws.Range("B2").Value = "Siddarth"
$ws.Range("C2").Value = "chandra.siddarth@yahoo.com"

# Update the active selection on the sheet from B7 to C7.
$ws.Range("C7").Select()

$wb.Save()
